$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Start row 10 off with the same look as the row above it (row 9), then
# fix up the couple of columns that need a different style and fill in
# the new data for the PPP2R5D / wrd entry.
$ws.Range("A9:K9").Copy() | Out-Null
$ws.Range("A10:K10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Columns E and G use the plain bordered look (same group as B/C/D) for
# this row rather than row 9's shaded look.
$ws.Range("B9").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("G10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# No "selected_fly_protein" value yet for this entry -- remove it entirely
# rather than leave an empty formatted cell.
$ws.Range("K10").Clear()

# Row height matches the rest of the table.
$ws.Rows.Item(10).RowHeight = 55

# --- New row of data for PPP2R5D / wrd (added at the bottom of the table) ---
$ws.Range("A10").Value = "PPP2R5D"
$ws.Range("B10").Value = "ENST00000485511"
$ws.Range("C10").Value = "NM_006245"
$ws.Range("D10").Value = "NP_006236"
$ws.Range("E10").Value = "CCDS4878"
$ws.Range("F10").Value = "wrd"
$ws.Range("G10").Value = "FBgn0042693"
$ws.Range("H10").Value = "many isoforms"
